$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Row 1 headers - add new process columns F1, G1
$ws.Range("F1").Value = "proMotorcycleInsurance"
$ws.Range("G1").Value = "proCamperInsurance"

# Row 3: Automobile smoke test (renamed from 102_VehicleInsuranceAutomobile_001_SmokeTest -> 102_AutomobileInsurance_001_SmokeTest)
$ws.Range("A3").Value = "102_AutomobileInsurance_001_SmokeTest"
$ws.Range("B3").Value = "var102_AutomobileInsurance_001_SmokeTest"
$ws.Range("D3").Value = "102_AutomobileInsurance_001_SmokeTest"

# Row 5: new Motorcycle smoke test
$ws.Range("A5").Value = "104_MotorcycleInsurance_001_SmokeTest"
$ws.Range("B5").Value = "var104_MotorcycleInsurance_001_SmokeTest"
$ws.Range("C5").Value = "Open Motorcycle Insurance"
$ws.Range("F5").Value = "104_MotorcycleInsurance_001_SmokeTest"

# Row 6: new Camper smoke test
$ws.Range("A6").Value = "105_CamperInsurance_001_SmokeTest"
$ws.Range("B6").Value = "var105_CamperInsurance_001_SmokeTest"
$ws.Range("C6").Value = "Open Camper Insurance"
$ws.Range("G6").Value = "105_CamperInsurance_001_SmokeTest"

# Autofit columns to reflect new content widths
$ws.Columns("A:G").AutoFit() | Out-Null

# Update selection to reflect final cursor position as in the diff
$ws.Range("F7").Select()
